$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13662
$ws1.Range("F4").Value = 125
$ws1.Range("F6").Value = 498
$ws1.Range("F9").Value = 13872
$ws1.Range("F10").Value = 14687
$ws1.Range("F20").Value = 14
$ws1.Range("F21").Value = 1139
$ws1.Range("F24").Value = 5664
$ws1.Range("F27").Value = 5387
$ws1.Range("F28").Value = 43
$ws1.Range("F29").Value = 43

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13662
$ws4.Range("F5").Value = 125
$ws4.Range("F7").Value = 498
$ws4.Range("F10").Value = 13872
$ws4.Range("F11").Value = 14687
$ws4.Range("F21").Value = 14
$ws4.Range("F22").Value = 1139
$ws4.Range("F25").Value = 5664
$ws4.Range("F28").Value = 5387
$ws4.Range("F29").Value = 43
$ws4.Range("F30").Value = 43
